# This script is driven against $excel.ActiveWorkbook (already open).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "user story 1": only the remembered cell-selection moves, D31 -> D32.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D32").Select()

# ---------------------------------------------------------------
# "user story 5": the stray trailing row (A16 = 4, with nothing else
# in the row) is removed, and the remembered selection moves to A25.
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A16").EntireRow.Delete()
$ws5.Range("A25").Select()

# ---------------------------------------------------------------
# New sheet "user story 6", appended after "user story 5". Built by
# duplicating "user story 5" (same layout/styles/merges/column widths)
# and then overwriting the cells that differ for this story.
# ---------------------------------------------------------------
$ws5.Copy([System.Reflection.Missing]::Value, $ws5)
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "user story 6"

$ws6.Range("B1").Value = "Issue #42"
$ws6.Range("B3").Value = "User story #6"
$ws6.Range("C3").Value = "As an user, I shall search one channel on the webpage."
$ws6.Range("B8").Value = "open one client"
$ws6.Range("B9").Value = "click search area"
$ws6.Range("B10").Value = "type in channel name"
$ws6.Range("C10").Value = "test1"
$ws6.Range("D10").Value = "The channel “test1” show up on the page"
$ws6.Range("E10").Value = "The channel “test1” show up on the page"
$ws6.Range("F10").Value = "✔"
$ws6.Range("B14").Value = "click search area"
$ws6.Range("B15").Value = "type in channel name"
$ws6.Range("C15").Value = "no exist"
$ws6.Range("D15").Value = "there is no channel show up"
$ws6.Range("E15").Value = "there is no channel show up"

# cells that exist in "user story 5" but have no counterpart in this story
$ws6.Range("D8:F8").Clear()
$ws6.Range("D9:F9").Clear()
$ws6.Range("B11").Clear()
$ws6.Range("D11:F11").Clear()
$ws6.Range("D13:F13").Clear()
$ws6.Range("D14:F14").Clear()
$ws6.Range("D28").Select()

# ---------------------------------------------------------------
# New sheet "user story 7", appended after "user story 6". Same
# duplication approach as "user story 6" above.
# ---------------------------------------------------------------
$ws5.Copy([System.Reflection.Missing]::Value, $ws6)
$ws7 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7.Name = "user story 7"

$ws7.Range("B1").Value = "Issue #75"
$ws7.Range("B3").Value = "User story #7"
$ws7.Range("C3").Value = "When typing a message apears at the bottom of the page indicating the possible actions:  Return to send message Shift + Return to add new line This message should disappear when input field is not in focus."
$ws7.Range("B8").Value = "open one client"
$ws7.Range("B9").Value = "click text area"
$ws7.Range("B10").Value = "focus text area"
$ws7.Range("D10").Value = "Tip show up"
$ws7.Range("E10").Value = "Tip show up"
$ws7.Range("F10").Value = "✔"
$ws7.Range("B13").Value = "open one client"
$ws7.Range("B14").Value = "click text area"
$ws7.Range("B15").Value = "unfocus text area"
$ws7.Range("D15").Value = "nothing will show up"
$ws7.Range("E15").Value = "nothing will show up"

# cells that exist in "user story 5" but have no counterpart in this story
$ws7.Range("D8:F8").Clear()
$ws7.Range("D9:F9").Clear()
$ws7.Range("C10").Clear()
$ws7.Range("B11").Clear()
$ws7.Range("D11:F11").Clear()
$ws7.Range("D13:F13").Clear()
$ws7.Range("D14:F14").Clear()
$ws7.Range("C15").Clear()
$ws7.Range("B4").Select()

# ---------------------------------------------------------------
# "user story 3" becomes the active tab (workbook activeTab = 2,
# i.e. the 3rd tab, 0-indexed). Must run last: adding/copying sheets
# above shifts the active tab to whichever sheet was just created.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
